$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Ensure column D retains text formatting for numeric-looking values
$colD = $ws.Range("D2:D51")
$colD.NumberFormat = "@"

$ws.Range('D2').Value = '26.956.03'
$ws.Range('D3').Value = '1.873.98'
$ws.Range('D5').Value = '305.38'
$ws.Range('E5').Value = '  -0.33%  '
$ws.Range('E6').Value = '  +0.02%  '
$ws.Range('D7').Value = '0.5089'
$ws.Range('E7').Value = '  -0.81%  '
$ws.Range('D8').Value = '0.3668'
$ws.Range('E8').Value = '  -2.26%  '
$ws.Range('D9').Value = '0.07207'
$ws.Range('E9').Value = '  +0.69%  '
$ws.Range('D10').Value = '0.8950'
$ws.Range('E10').Value = '  +0.75%  '
$ws.Range('D11').Value = '20.76'
$ws.Range('E11').Value = '  +0.43%  '
$ws.Range('D12').Value = '1.886.01'
$ws.Range('E12').Value = '  +1.05%  '
$ws.Range('D13').Value = '0.07524'
$ws.Range('E13').Value = '  -1.01%  '
$ws.Range('D14').Value = '95.03'
$ws.Range('E14').Value = '  +6.15%  '
$ws.Range('D15').Value = '5.245'
$ws.Range('E15').Value = '  -1.09%  '
$ws.Range('D16').Value = '1.001'
$ws.Range('D17').Value = '0.000008524'
$ws.Range('E17').Value = '  +0.63%  '
$ws.Range('D18').Value = '14.25'
$ws.Range('E18').Value = '  +1.24%  '
$ws.Range('D19').Value = '0.9998'
$ws.Range('E19').Value = '  -0.08%  '
$ws.Range('D20').Value = '27.002.43'
$ws.Range('E20').Value = '  -0.24%  '
$ws.Range('D21').Value = '5.024'
$ws.Range('E21').Value = '  -0.14%  '
$ws.Range('D22').Value = '2.135.33'
$ws.Range('E22').Value = '  +2.67%  '
$ws.Range('E23').Value = '  -1.17%  '
$ws.Range('D24').Value = '6.399'
$ws.Range('E24').Value = '  -0.93%  '
$ws.Range('D25').Value = '148.48'
$ws.Range('E25').Value = '  +0.75%  '
$ws.Range('D26').Value = '1.794'
$ws.Range('E26').Value = '  -2.63%  '
$ws.Range('E27').Value = '  -0.39%  '
$ws.Range('D28').Value = '2.089'
$ws.Range('E28').Value = '  -1.15%  '
$ws.Range('D29').Value = '113.54'
$ws.Range('E29').Value = '  +0.68%  '
$ws.Range('D30').Value = '4.728'
$ws.Range('E30').Value = '  +1.41%  '
$ws.Range('D31').Value = '4.729'
$ws.Range('E31').Value = '  +0.72%  '
$ws.Range('E32').Value = '  +0.66%  '
$ws.Range('D33').Value = '0.05087'
$ws.Range('E33').Value = '  -0.83%  '
$ws.Range('D34').Value = '0.7512'
$ws.Range('E34').Value = '  +3.50%  '
$ws.Range('D35').Value = '2.965'
$ws.Range('E35').Value = '  -3.25%  '
$ws.Range('D36').Value = '1.159'
$ws.Range('E36').Value = '  +0.25%  '
$ws.Range('D37').Value = '3.236'
$ws.Range('E37').Value = '  +6.37%  '
$ws.Range('D38').Value = '2.531'
$ws.Range('E38').Value = '  +1.64%  '
$ws.Range('D39').Value = '0.5638'
$ws.Range('E39').Value = '  +5.68%  '
$ws.Range('E40').Value = '  -1.76%  '
$ws.Range('E41').Value = '  +0.45%  '
$ws.Range('D42').Value = '6.633'
$ws.Range('E42').Value = '  +1.15%  '
$ws.Range('D43').Value = '115.44'
$ws.Range('E43').Value = '  -1.44%  '
$ws.Range('D44').Value = '8.582'
$ws.Range('E44').Value = '  +3.57%  '
$ws.Range('E45').Value = '  +0.77%  '
$ws.Range('D46').Value = '0.4753'
$ws.Range('E46').Value = '  +2.42%  '
$ws.Range('B47').Value = 'PaxDollar'
$ws.Range('C47').Value = 'https://coinranking.com/coin/JCKLgWPAF+paxdollar-usdp'
$ws.Range('D47').Value = '0.9997'
$ws.Range('E47').Value = '  +0.01%  '
$ws.Range('B48').Value = 'EnergySwap'
$ws.Range('C48').Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
$ws.Range('D48').Value = '10.12'
$ws.Range('E48').Value = '  +1.38%  '
$ws.Range('D49').Value = '1.572'
$ws.Range('E49').Value = '  +0.09%  '
$ws.Range('D50').Value = '36.98'
$ws.Range('E50').Value = '  +1.19%  '
$ws.Range('D51').Value = '63.22'
$ws.Range('E51').Value = '  -1.11%  '

# Restore default style on column D so no stray number-format remains
$colD.Style = "Normal"

Write-Host "Done applying crypto list updates"